$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 230
$ws.Range("I6").Value = 230
$ws.Range("K6").Value = 690
$ws.Range("M6").Value = -578
$ws.Range("H112").Value = 2825.2
$ws.Range("J112").Value = 3050.2222
$ws.Range("L112").Value = 9150.6666
$ws.Range("N112").Value = -11366.6666
$ws.Range("H113").Value = 3600
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3600
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3600
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10108
$ws.Range("H115").Value = 2211.3333
$ws.Range("I115").Value = 827
$ws.Range("K115").Value = 2481
$ws.Range("M115").Value = -914
$ws.Range("H127").Value = 942.449
$ws.Range("I127").Value = 231.85715
$ws.Range("K127").Value = 695.5714499999999
$ws.Range("M127").Value = 4264.428550000001
$ws.Range("H129").Value = 936.2542
$ws.Range("I129").Value = 484
$ws.Range("J129").Value = 1076.9556
$ws.Range("K129").Value = 1452
$ws.Range("L129").Value = 3230.8668
$ws.Range("M129").Value = 3548
$ws.Range("N129").Value = -13230.8668
$ws.Range("H132").Value = 3056.762
$ws.Range("I132").Value = 2636.6155
$ws.Range("J132").Value = 3739.5
$ws.Range("K132").Value = 7909.8465
$ws.Range("L132").Value = 11218.5
$ws.Range("M132").Value = -5379.8465
$ws.Range("N132").Value = -16278.5
$ws.Range("H138").Value = 3389.7856
$ws.Range("I138").Value = 2103.1428
$ws.Range("J138").Value = 3711.4465
$ws.Range("K138").Value = 6309.428400000001
$ws.Range("L138").Value = 11134.3395
$ws.Range("M138").Value = -1169.428400000001
$ws.Range("N138").Value = -21414.3395

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 30118.334
$ws.Range("J24").Value = 30118.334
$ws.Range("L24").Value = 30118.334
$ws.Range("N24").Value = -30866.334
$ws.Range("H32").Value = 9798.385
$ws.Range("I32").Value = 10127.339
$ws.Range("K32").Value = 10127.339
$ws.Range("M32").Value = -9840.339
$ws.Range("H100").Value = 30118.334
$ws.Range("J100").Value = 30118.334
$ws.Range("L100").Value = 30118.334
$ws.Range("N100").Value = -32282.334
$ws.Range("H122").Value = 2750
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1318
$ws.Range("I16").Value = 1429.6666
$ws.Range("J16").Value = 1150.5
$ws.Range("K16").Value = 1429.6666
$ws.Range("L16").Value = 1150.5
$ws.Range("M16").Value = -1142.6666
$ws.Range("N16").Value = -1724.5
$ws.Range("H31").Value = 2103.0159
$ws.Range("I31").Value = 1732.5714
$ws.Range("J31").Value = 3399.5715
$ws.Range("K31").Value = 1732.5714
$ws.Range("L31").Value = 3399.5715
$ws.Range("M31").Value = -1437.5714
$ws.Range("N31").Value = -3989.5715
$ws.Range("H34").Value = 2103.0159
$ws.Range("I34").Value = 1732.5714
$ws.Range("J34").Value = 3399.5715
$ws.Range("K34").Value = 1732.5714
$ws.Range("L34").Value = 3399.5715
$ws.Range("M34").Value = -1530.5714
$ws.Range("N34").Value = -3803.5715
$ws.Range("H99").Value = 4330
$ws.Range("I99").Value = 5000
$ws.Range("J99").Value = 3660
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 3660
$ws.Range("M99").Value = -3502
$ws.Range("N99").Value = -6656
$ws.Range("H113").Value = 1318
$ws.Range("I113").Value = 1429.6666
$ws.Range("J113").Value = 1150.5
$ws.Range("K113").Value = 1429.6666
$ws.Range("L113").Value = 1150.5
$ws.Range("M113").Value = 740.3334
$ws.Range("N113").Value = -5490.5
$ws.Range("H126").Value = 4330
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 3660
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 10980
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -15920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1601.0385
$ws.Range("J5").Value = 1038.3077
$ws.Range("L5").Value = 3114.9231
$ws.Range("N5").Value = -3338.9231
$ws.Range("H38").Value = 236.66667
$ws.Range("I38").Value = 224.54546
$ws.Range("J38").Value = 290
$ws.Range("K38").Value = 673.6363799999999
$ws.Range("L38").Value = 870
$ws.Range("M38").Value = -326.6363799999999
$ws.Range("N38").Value = -1564
$ws.Range("H103").Value = 2229
$ws.Range("I103").Value = 656
$ws.Range("J103").Value = 5375
$ws.Range("K103").Value = 1968
$ws.Range("L103").Value = 16125
$ws.Range("M103").Value = -1089
$ws.Range("N103").Value = -17883
$ws.Range("H121").Value = 869.0769
$ws.Range("I121").Value = 592.375
$ws.Range("J121").Value = 940.4838999999999
$ws.Range("K121").Value = 1777.125
$ws.Range("L121").Value = 2821.4517
$ws.Range("M121").Value = -467.125
$ws.Range("N121").Value = -5441.4517
$ws.Range("H131").Value = 21740770
$ws.Range("J131").Value = 23257532
$ws.Range("L131").Value = 69772596
$ws.Range("N131").Value = -69782676
$ws.Range("H135").Value = 1601.0385
$ws.Range("J135").Value = 1038.3077
$ws.Range("L135").Value = 9344.7693
$ws.Range("N135").Value = -14414.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 14000
$ws.Range("I19").Value = 20000
$ws.Range("J19").Value = 8000
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 8000
$ws.Range("M19").Value = -19712
$ws.Range("N19").Value = -8576
$ws.Range("H122").Value = 2773.5881
$ws.Range("I122").Value = 3975.7778
$ws.Range("K122").Value = 11927.3334
$ws.Range("M122").Value = -9477.3334
$ws.Range("H126").Value = 3195.1333
$ws.Range("J126").Value = 3628.5
$ws.Range("L126").Value = 10885.5
$ws.Range("N126").Value = -15825.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4036.5386
$ws.Range("J7").Value = 3725
$ws.Range("L7").Value = 3725
$ws.Range("N7").Value = -3949
$ws.Range("H23").Value = 1500
$ws.Range("I23").Value = 1500
$ws.Range("K23").Value = 1500
$ws.Range("M23").Value = -1270
$ws.Range("H40").Value = 17000
$ws.Range("I40").Value = 30000
$ws.Range("K40").Value = 30000
$ws.Range("M40").Value = -29864
$ws.Range("H101").Value = 117720
$ws.Range("J101").Value = 117720
$ws.Range("L101").Value = 117720
$ws.Range("N101").Value = -124210
$ws.Range("H126").Value = 4036.5386
$ws.Range("J126").Value = 3725
$ws.Range("L126").Value = 11175
$ws.Range("N126").Value = -16115

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 27500
$ws.Range("H81").Value = 169866.67
$ws.Range("I81").Value = 202840
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 405680
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -404619
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 169866.67
$ws.Range("I84").Value = 202840
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 2028400
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -2023096
$ws.Range("N84").Value = -60608
$ws.Range("H95").Value = 200000
$ws.Range("J95").Value = 200000
$ws.Range("L95").Value = 200000
$ws.Range("N95").Value = -205492
$ws.Range("H126").Value = 13655.929
$ws.Range("I126").Value = 18786.3
$ws.Range("J126").Value = 830
$ws.Range("K126").Value = 56358.89999999999
$ws.Range("L126").Value = 2490
$ws.Range("M126").Value = -53888.89999999999
$ws.Range("N126").Value = -7430
$ws.Range("H136").Value = 956.19354
$ws.Range("I136").Value = 912.7037
$ws.Range("J136").Value = 1249.75
$ws.Range("K136").Value = 2738.1111
$ws.Range("L136").Value = 3749.25
$ws.Range("M136").Value = -188.1111000000001
$ws.Range("N136").Value = -8849.25
